$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "609-3273-ND"
$ws.Range("E6").Value = "FCI"
$ws.Range("F6").Value = "68001-108HLF"
$ws.Range("I6").Value = "1x8"
$ws.Range("C6").Value = "1x8 Headers Gold 30uin"
$ws.Range("H6").Value = "0.230in"
$ws.Range("D8").Select() | Out-Null
